$d = $word.ActiveDocument

# The opening paragraph ("This Mutual Non-Disclosure Agreement (this
# "Agreement") is entered into ... {{ party2_address }}.") needs several of
# its runs re-split so that:
#   - "between" is split across two runs (mid-word split, cosmetic only)
#   - the party1/party2 placeholder runs are split into a "name" run (kept
#     bold) and a separate "conditional entity / address" run (not bold)
#   - the ", and " separator becomes two runs: ", and" then " "
#
# Find.Execute/InsertXML on a sub-range of a paragraph corrupts neighboring
# run text in this runtime, so instead we rebuild the *entire* paragraph's
# run content (everything except the trailing paragraph mark) in one shot.
$p = $d.Paragraphs(2)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)

$expectedText = "This Mutual Non-Disclosure Agreement (this " + [char]0x201C + "Agreement" + [char]0x201D + ") is entered into and made effective as of __________, 20__, between {{ party1_name }}{% if party1_entity != " + [char]34 + "Individual" + [char]34 + " %}, a {{ party1_state }} {{ party1_entity }}{% endif %}, whose address is {{ party1_address }}, and {{ party2_name }}{% if party2_entity != " + [char]34 + "Individual" + [char]34 + " %}, a {{ party2_state }} {{ party2_entity }}{% endif %}, whose address is {{ party2_address }}."
if ($rng.Text -ne $expectedText) {
    throw "Paragraph 2 text did not match the expected NDA intro sentence; aborting to avoid corrupting the wrong paragraph."
}

$boldPartyRPr = '<w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue"/><w:b/><w:bCs/><w:color w:val="0E0E0E"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>'
$plainPartyRPr = '<w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue"/><w:color w:val="0E0E0E"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>'

$runs = ''
$runs += '<w:r><w:t>This Mutual Non-Disclosure Agreement (this &#8220;</w:t></w:r>'
$runs += '<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Agreement</w:t></w:r>'
$runs += '<w:r><w:t>&#8221;) is entered into and made effective as of __________, 20__, b</w:t></w:r>'
$runs += '<w:r><w:t xml:space="preserve">etween </w:t></w:r>'
$runs += '<w:r>' + $boldPartyRPr + '<w:t>{{ party1_name }}</w:t></w:r>'
$runs += '<w:r>' + $plainPartyRPr + '<w:t>{% if party1_entity != &quot;Individual&quot; %}, a {{ party1_state }} {{ party1_entity }}{% endif %}, whose address is {{ party1_address }}</w:t></w:r>'
$runs += '<w:r><w:t>, and</w:t></w:r>'
$runs += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$runs += '<w:r>' + $boldPartyRPr + '<w:t>{{ party2_name }}</w:t></w:r>'
$runs += '<w:r>' + $plainPartyRPr + '<w:t>{% if party2_entity != &quot;Individual&quot; %}, a {{ party2_state }} {{ party2_entity }}{% endif %}, whose address is {{ party2_address }}</w:t></w:r>'
$runs += '<w:r><w:t>.</w:t></w:r>'

$xmlPayload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $runs + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xmlPayload)
